$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 13890895
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 14707948
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 44123844
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -44124384
$ws.Range("H73").Value = 13890895
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 14707948
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 44123844
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -44125716
$ws.Range("H76").Value = 4999.8
$ws.Range("I76").Value = 4999.857
$ws.Range("J76").Value = 4999
$ws.Range("K76").Value = 4999.857
$ws.Range("L76").Value = 4999
$ws.Range("M76").Value = -4684.857
$ws.Range("N76").Value = -5629
$ws.Range("H79").Value = 4999.8
$ws.Range("I79").Value = 4999.857
$ws.Range("J79").Value = 4999
$ws.Range("K79").Value = 4999.857
$ws.Range("L79").Value = 4999
$ws.Range("M79").Value = -3907.857
$ws.Range("N79").Value = -7183
$ws.Range("H86").Value = 140627220
$ws.Range("I86").Value = 76925600
$ws.Range("J86").Value = 416667500
$ws.Range("K86").Value = 76925600
$ws.Range("L86").Value = 416667500
$ws.Range("M86").Value = -76924477
$ws.Range("N86").Value = -416669746
$ws.Range("H89").Value = 140627220
$ws.Range("I89").Value = 76925600
$ws.Range("J89").Value = 416667500
$ws.Range("K89").Value = 384628000
$ws.Range("L89").Value = 2083337500
$ws.Range("M89").Value = -384622384
$ws.Range("N89").Value = -2083348732
$ws.Range("H129").Value = 2043.3334
$ws.Range("I129").Value = 1000
$ws.Range("J129").Value = 2173.75
$ws.Range("K129").Value = 3000
$ws.Range("L129").Value = 6521.25
$ws.Range("M129").Value = 2000
$ws.Range("N129").Value = -16521.25
$ws.Range("H132").Value = 2997.2034
$ws.Range("I132").Value = 2997.2034
$ws.Range("K132").Value = 8991.610199999999
$ws.Range("M132").Value = -6461.610199999999
$ws.Range("H137").Value = 21909.791
$ws.Range("J137").Value = 3076.077
$ws.Range("L137").Value = 9228.231
$ws.Range("N137").Value = -14328.231
$ws.Range("H138").Value = 2859.24
$ws.Range("I138").Value = 1315.8966
$ws.Range("J138").Value = 3489.6196
$ws.Range("K138").Value = 3947.6898
$ws.Range("L138").Value = 10468.8588
$ws.Range("M138").Value = 1192.3102
$ws.Range("N138").Value = -20748.8588
$ws.Range("H141").Value = 32308.691
$ws.Range("I141").Value = 1637.5454
$ws.Range("J141").Value = 201000
$ws.Range("K141").Value = 4912.6362
$ws.Range("L141").Value = 603000
$ws.Range("M141").Value = 267.3638000000001
$ws.Range("N141").Value = -613360

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 333
$ws.Range("I5").Value = 340.3846
$ws.Range("K5").Value = 340.3846
$ws.Range("M5").Value = -228.3846
$ws.Range("H32").Value = 26986562
$ws.Range("I32").Value = 27092512
$ws.Range("J32").Value = 25644524
$ws.Range("K32").Value = 27092512
$ws.Range("L32").Value = 25644524
$ws.Range("M32").Value = -27092225
$ws.Range("N32").Value = -25645098
$ws.Range("H45").Value = 4527.273
$ws.Range("I45").Value = 4637.5
$ws.Range("K45").Value = 4637.5
$ws.Range("M45").Value = -4260.5
$ws.Range("H61").Value = 1847.9231
$ws.Range("I61").Value = 1623.84
$ws.Range("K61").Value = 1623.84
$ws.Range("M61").Value = -1411.84
$ws.Range("H97").Value = 726.4643
$ws.Range("I97").Value = 579.0909
$ws.Range("K97").Value = 579.0909
$ws.Range("M97").Value = -83.09090000000003
$ws.Range("H132").Value = 143179.58
$ws.Range("I132").Value = 174300.53
$ws.Range("J132").Value = 4332.231
$ws.Range("K132").Value = 522901.59
$ws.Range("L132").Value = 12996.693
$ws.Range("M132").Value = -520371.59
$ws.Range("N132").Value = -18056.693
$ws.Range("H136").Value = 1847.9231
$ws.Range("I136").Value = 1623.84
$ws.Range("K136").Value = 4871.52
$ws.Range("M136").Value = -2321.52

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 333
$ws.Range("I4").Value = 340.3846
$ws.Range("K4").Value = 340.3846
$ws.Range("M4").Value = -225.3846
$ws.Range("H107").Value = 1516.7916
$ws.Range("I107").Value = 863.3158
$ws.Range("K107").Value = 863.3158
$ws.Range("M107").Value = 1056.6842
$ws.Range("H134").Value = 2103134.8
$ws.Range("I134").Value = 2383019.8
$ws.Range("J134").Value = 3997.5
$ws.Range("K134").Value = 7149059.399999999
$ws.Range("L134").Value = 11992.5
$ws.Range("M134").Value = -7146524.399999999
$ws.Range("N134").Value = -17062.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8039.2
$ws.Range("I31").Value = 5799
$ws.Range("J31").Value = 17000
$ws.Range("K31").Value = 5799
$ws.Range("L31").Value = 17000
$ws.Range("M31").Value = -5504
$ws.Range("N31").Value = -17590
$ws.Range("H34").Value = 8039.2
$ws.Range("I34").Value = 5799
$ws.Range("J34").Value = 17000
$ws.Range("K34").Value = 5799
$ws.Range("L34").Value = 17000
$ws.Range("M34").Value = -5597
$ws.Range("N34").Value = -17404
$ws.Range("H58").Value = 2344.9565
$ws.Range("J58").Value = 3186.6
$ws.Range("L58").Value = 3186.6
$ws.Range("N58").Value = -3592.6
$ws.Range("H62").Value = 2628.4285
$ws.Range("I62").Value = 2474.75
$ws.Range("K62").Value = 2474.75
$ws.Range("M62").Value = -1850.75
$ws.Range("H65").Value = 2628.4285
$ws.Range("I65").Value = 2474.75
$ws.Range("K65").Value = 12373.75
$ws.Range("M65").Value = -9253.75
$ws.Range("H105").Value = 1927.6923
$ws.Range("I105").Value = 1796.6666
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 1796.6666
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -49.66660000000002
$ws.Range("N105").Value = -6994
$ws.Range("H122").Value = 4005098.5
$ws.Range("H132").Value = 3972.8147
$ws.Range("I132").Value = 3763.6978
$ws.Range("K132").Value = 11291.0934
$ws.Range("M132").Value = -8761.0934
$ws.Range("H134").Value = 2441.0527
$ws.Range("I134").Value = 2390.361
$ws.Range("J134").Value = 2527.9524
$ws.Range("K134").Value = 7171.083
$ws.Range("L134").Value = 7583.8572
$ws.Range("M134").Value = -4636.083
$ws.Range("N134").Value = -12653.8572
$ws.Range("H136").Value = 2344.9565
$ws.Range("J136").Value = 3186.6
$ws.Range("L136").Value = 9559.799999999999
$ws.Range("N136").Value = -14659.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 88129740
$ws.Range("I4").Value = 83835336
$ws.Range("K4").Value = 251506008
$ws.Range("M4").Value = -251505896

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1950.2
$ws.Range("I80").Value = 1884.8334
$ws.Range("J80").Value = 2048.25
$ws.Range("K80").Value = 1884.8334
$ws.Range("L80").Value = 2048.25
$ws.Range("M80").Value = -886.8334
$ws.Range("N80").Value = -4044.25
$ws.Range("H83").Value = 1950.2
$ws.Range("I83").Value = 1884.8334
$ws.Range("J83").Value = 2048.25
$ws.Range("K83").Value = 9424.166999999999
$ws.Range("L83").Value = 10241.25
$ws.Range("M83").Value = -4432.166999999999
$ws.Range("N83").Value = -20225.25
$ws.Range("H102").Value = 2051.8333
$ws.Range("I102").Value = 1956.5555
$ws.Range("J102").Value = 2337.6667
$ws.Range("K102").Value = 1956.5555
$ws.Range("L102").Value = 2337.6667
$ws.Range("M102").Value = -334.5554999999999
$ws.Range("N102").Value = -5581.6667
$ws.Range("H122").Value = 1777.8
$ws.Range("I122").Value = 1296.5555
$ws.Range("J122").Value = 2499.6667
$ws.Range("K122").Value = 3889.6665
$ws.Range("L122").Value = 7499.000100000001
$ws.Range("M122").Value = -1439.6665
$ws.Range("N122").Value = -12399.0001
$ws.Range("H125").Value = 99232.25
$ws.Range("J125").Value = 99232.25
$ws.Range("L125").Value = 99232.25
$ws.Range("N125").Value = -104152.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5582.8535
$ws.Range("I46").Value = 6595.174
$ws.Range("K46").Value = 6595.174
$ws.Range("M46").Value = -6407.174
$ws.Range("H82").Value = 12477.105
$ws.Range("I82").Value = 1280.2222
$ws.Range("J82").Value = 22554.3
$ws.Range("K82").Value = 1280.2222
$ws.Range("L82").Value = 22554.3
$ws.Range("M82").Value = -919.2221999999999
$ws.Range("N82").Value = -23276.3
$ws.Range("H85").Value = 12477.105
$ws.Range("I85").Value = 1280.2222
$ws.Range("J85").Value = 22554.3
$ws.Range("K85").Value = 1280.2222
$ws.Range("L85").Value = 22554.3
$ws.Range("M85").Value = -32.22219999999993
$ws.Range("N85").Value = -25050.3
$ws.Range("H132").Value = 385532.06
$ws.Range("I132").Value = 422196.47
$ws.Range("K132").Value = 1266589.41
$ws.Range("M132").Value = -1264059.41

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 34499
$ws.Range("J69").Value = 34499
$ws.Range("L69").Value = 34499
$ws.Range("N69").Value = -35997
$ws.Range("H72").Value = 34499
$ws.Range("J72").Value = 34499
$ws.Range("L72").Value = 103497
$ws.Range("N72").Value = -110985
$ws.Range("H122").Value = 27780702
$ws.Range("I122").Value = 31252150
$ws.Range("K122").Value = 93756450
$ws.Range("M122").Value = -93754000
$ws.Range("H132").Value = 21182.527
$ws.Range("I132").Value = 26949.65
$ws.Range("K132").Value = 80848.95000000001
$ws.Range("M132").Value = -78318.95000000001
$ws.Range("H136").Value = 17382.127
$ws.Range("J136").Value = 69599.734
$ws.Range("L136").Value = 208799.202
$ws.Range("N136").Value = -213899.202
